$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.604.38'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.193.82'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '260.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '82.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +12.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.11%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.592'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('E10').Value = '  +9.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0918'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.96'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.56%  '
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.518.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.175.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.772'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.537.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.15%  '
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +16.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.45'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.86'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.20%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '42.06'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +16.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.06%  '
$ws.Range('B28').Value = 'WEMIXToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.30%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.42'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.55%  '
$ws.Range('E33').Value = '  +7.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.34'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.17%  '
$ws.Range('E35').Value = '  +7.96%  '
$ws.Range('E36').Value = '  +2.41%  '
$ws.Range('E37').Value = '  +7.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0352'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.11'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +17.89%  '
$ws.Range('E41').Value = '  +2.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '63.74'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.45'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.200'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.09%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0982'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.37'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('E48').Value = '  +4.31%  '
$ws.Range('E49').Value = '  +3.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.437'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.03%  '
$ws.Range('E51').Value = '  +24.16%  '
